$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("begroting")
$ws1.Range("A21").Value = "hours"
$ws1.Range("B21").Formula = "=24480/135"
$ws1.Range("A22").Value = "booked"
$ws1.Range("B22").Formula = "=8+34"
$ws1.Range("A23").Value = "remaining"
$ws1.Range("B23").Formula = "=B21-B22"
$ws1.Range("C23").Formula = "=(B21-B22)/B21*100"
$ws1.Range("C23").NumberFormat = "0.0"
$ws2 = $wb.Worksheets.Item("api calls")
$ws2.Range("B12").Select()
$ws1.Activate()
$ws1.Range("C24").Select()

